# ----------------------------------------------------------------------
# Apply the "Update tanggal 22 Mei 2019" revisions to BAB I.docx
#
#   * Drop the "Menurut " lead-in on four in-text citations (Trinova,
#     Khobir, Novaliendry, Primasari) and instead attach an explicit
#     "menyatakan bahwa" (",") clause after the cited work, matching the
#     surrounding prose style used elsewhere in the chapter.
#   * Terminate the "(Ekawati & Falani, 2015)" sentence with a period.
#   * Move the Primasari citation to a parenthetical "(Primasari , 2013)"
#     right after "... sistem operasi android" and capitalise the
#     paragraph's new opening word ("Indonesia").
#   * The document's "last edit" bookmark (_GoBack) therefore now sits at
#     that citation instead of on the "Gambar 7" caption.
# ----------------------------------------------------------------------
$d = $word.ActiveDocument

function FindReplace($findText, $replaceText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Output "WARNING: replace failed for [$findText]"
    }
    return $ok
}

# 1) Trinova paragraph: drop leading "Menurut " and add "menyatakan bahwa"
#    before the trailing ", guru ..." clause.
FindReplace "Menurut Trinova" "Trinova" | Out-Null
FindReplace "peserta didik, guru dapat menerapkan" "peserta didik menyatakan bahwa, guru dapat menerapkan" | Out-Null

# 2) Close the "(Ekawati & Falani, 2015)" paragraph with a period.
FindReplace "(Ekawati & Falani, 2015)" "(Ekawati & Falani, 2015)." | Out-Null

# 3) Khobir paragraph: drop leading "Menurut " and add "menyatakan bahwa"
#    before the ", permainan edukatif ..." clause.
FindReplace "Menurut Khobir (2009)" "Khobir (2009)" | Out-Null
FindReplace "permainan edukatif, permainan edukatif hendaknya" "permainan edukatif menyatakan bahwa, permainan edukatif hendaknya" | Out-Null

# 4) Novaliendry paragraph: drop leading "Menurut " and add
#    "menyatakan bahwa," before "game edukasi".
FindReplace "Menurut Novaliendry (2013)" "Novaliendry (2013) menyatakan bahwa," | Out-Null

# 5) Primasari paragraph: drop leading "Menurut " + lower-case "i" becomes
#    capital "I".
FindReplace "Menurut Primasari (2013) indonesia" "Indonesia" | Out-Null

#    ... and append the "(Primasari , 2013)" citation right after the
#    italic "android" run, without disturbing its formatting.
$androidRng = $d.Content
$androidOk = $androidRng.Find.Execute("telepon genggam dengan sistem operasi android", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $androidOk) {
    Write-Output "WARNING: could not locate the android run"
}
$androidRng.Collapse(0)
$citeStart = $androidRng.Start
$androidRng.InsertAfter(" (Primasari , 2013)")
$androidRng.Font.Italic = 0

# The last edit position ("_GoBack") now belongs right before "(Primasari"
# (i.e. one character -- the space -- after the end of "android"), so move
# the bookmark there from its old spot on the "Gambar 7" caption. Adding a
# bookmark with a name that already exists relocates it (Word only allows
# one bookmark per name), so this also removes it from "Gambar 7".
$bmRng = $d.Range($citeStart + 1, $citeStart + 1)
$d.Bookmarks.Add("_GoBack", $bmRng) | Out-Null

# 6) Touch the Gambar 4/5/6 captions so Word folds the adjacent plain runs
#    back together (same text, no content change) -- mirrors the caption
#    clean-up in the source revision.
FindReplace "Gambar 4" "Gambar 4" | Out-Null
FindReplace "Gambar 5" "Gambar 5" | Out-Null
FindReplace "Gambar 6" "Gambar 6" | Out-Null

Write-Output "done"
